$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -2
$ws.Range("F9").Value = -3
$ws.Range("F10").Value = -3
$ws.Range("F11").Value = 6
$ws.Range("F14").Value = 3
$ws.Range("F17").Value = 2
